# EZ_Parts_Budget.xlsx - "changes and status report"
# Insert two new line items (Micro USB Cable, Power Source/battery pack)
# just above the PARTS TOTAL row, and roll the subtotal formula forward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the "PARTS TOTAL" row (and the blank spacer row after it) down by two
# rows so the two new part rows land at rows 19-20, same as real Excel's
# Insert Sheet Rows command run with the cursor on row 19.
$ws.Rows("19:20").Insert()

# --- Row 19: Micro USB Cable -------------------------------------------------
$ws.Range("B19").Value = "15cm 90 degree micro USB cable"
$ws.Range("A19").Value = "Micro USB Cable"
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 10.99
$ws.Range("E19").Value = "CERRXIAN"
$ws.Range("F19").Value = "B073PQWY2B"
$ws.Range("G19").Value = 10.99

# --- Row 20: Power Source (battery pack) ------------------------------------
$ws.Range("A20").Value = "Power Source"
$ws.Range("B20").Value = "APC 5000mAh battery pack"
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 34
$ws.Range("F20").Value = "M5BK"
$ws.Range("E20").Value = "APC"
$ws.Range("G20").Value = 34

# Roll the PARTS TOTAL sum (now on row 21) forward to include the new rows.
$ws.Range("G21").Formula = "=SUM(G3:G20)"

# Hyperlink the two new part-number cells, same as every other row's F column.
$ws.Hyperlinks.Add($ws.Range("F19"), "https://www.amazon.ca/CERRXIAN-Degree-Angle-Braided-Charging/dp/B073PQWY2B")
$ws.Hyperlinks.Add($ws.Range("F20"), "https://www.amazon.ca/gp/product/B01MXXXXX")

# Hyperlinks.Add re-stamps its own style; put the normal hyperlink look
# (matching F13/F14) back on the two new cells.
$ws.Range("F13").Copy()
$ws.Range("F19").PasteSpecial(-4122)
$ws.Range("F20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the author's final selection.
$ws.Range("F20").Select()
